$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Row for year 2,008 (vars = 1)
Replace-Exact "3,111" "3,112"
Replace-Exact "0.4161304" "0.4161188"
Replace-Exact "0.1439181" "0.1438965"
Replace-Exact "0.4106085" "0.4105962"
Replace-Exact "1.4051915" "1.4056119"
Replace-Exact "16.0274895" "16.0332938"
Replace-Exact "0.002580272" "0.002579469"

# Row for year 2,012 (vars = 2)
Replace-Exact "3,111" "3,112"
Replace-Exact "0.3846426" "0.3848180"
Replace-Exact "0.1473129" "0.1476132"
Replace-Exact "0.3775934" "0.3776749"
Replace-Exact "0.91" "0.93"
Replace-Exact "0.88" "0.90"
Replace-Exact "0.4853522" "0.4948805"
Replace-Exact "0.1056934" "0.1370057"
Replace-Exact "0.002641135" "0.002646096"

# Row for year 2,016 (vars = 3)
Replace-Exact "3,111" "3,112"
Replace-Exact "0.3168062" "0.3169810"
Replace-Exact "0.1528448" "0.1531301"
Replace-Exact "0.3012580" "0.3013541"
Replace-Exact "0.9331339" "0.9385506"
Replace-Exact "0.6124881" "0.6310737"
Replace-Exact "0.002740316" "0.002744991"

Write-Output "Done applying replacements"
